$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 5.368399999999999
$ws.Range("B3").Value = 5.71839999999999
$ws.Range("E4").Value = 12.6931
$ws.Range("B5").Value = 5.2195
$ws.Range("E6").Value = 12.1522
$ws.Range("C7").Value = -11.49149999999999
$ws.Range("A9").Value = -20.40319999999998
$ws.Range("C9").Value = -12.33190000000001
$ws.Range("E10").Value = 11.8099
$ws.Range("B11").Value = 4.880100000000002
$ws.Range("E11").Value = 12.56749999999999
$ws.Range("B12").Value = 5.434599999999997
$ws.Range("A13").Value = -22.13510000000002
$ws.Range("A16").Value = -20.0333
$ws.Range("A18").Value = -22.65850000000001
$ws.Range("A20").Value = -22.08480000000003
$ws.Range("B21").Value = 5.074300000000003
$ws.Range("C21").Value = -13.7628
$ws.Range("E21").Value = 12.80729999999998
$ws.Range("E25").Value = 12.8488
